$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = 44967
$ws.Range("J2").Value = 400
$ws.Range("K2").Value = 22000
$ws.Range("L2").Value = 23000
$ws.Range("M2").Value = 22500
$ws.Range("P2").Value = 1731
$ws.Range("D3").Value = 44658
$ws.Range("J3").Value = 400
$ws.Range("K3").Value = 15000
$ws.Range("L3").Value = 16000
$ws.Range("M3").Value = 15500
$ws.Range("P3").Value = 1192
$ws.Range("D4").Value = 44435
$ws.Range("J4").Value = 480
$ws.Range("K4").Value = 13000
$ws.Range("L4").Value = 14000
$ws.Range("M4").Value = 13500
$ws.Range("P4").Value = 1038
$ws.Range("D5").Value = 44498
$ws.Range("J5").Value = 400
$ws.Range("K5").Value = 14000
$ws.Range("L5").Value = 15000
$ws.Range("M5").Value = 14500
$ws.Range("P5").Value = 1115
$ws.Range("D6").Value = 45044
$ws.Range("J6").Value = 400
$ws.Range("K6").Value = 16500
$ws.Range("L6").Value = 17000
$ws.Range("M6").Value = 16750
$ws.Range("P6").Value = 1288
$ws.Range("D7").Value = 44410
$ws.Range("J7").Value = 600
$ws.Range("K7").Value = 14000
$ws.Range("L7").Value = 15000
$ws.Range("M7").Value = 14500
$ws.Range("P7").Value = 1115
$ws.Range("D8").Value = 44874
$ws.Range("J8").Value = 520
$ws.Range("K8").Value = 13500
$ws.Range("L8").Value = 14000
$ws.Range("M8").Value = 13750
$ws.Range("P8").Value = 1058
$ws.Range("D9").Value = 44974
$ws.Range("J9").Value = 300
$ws.Range("K9").Value = 20000
$ws.Range("L9").Value = 21000
$ws.Range("M9").Value = 20500
$ws.Range("P9").Value = 1577
$ws.Range("D10").Value = 44806
$ws.Range("J10").Value = 500
$ws.Range("K10").Value = 14000
$ws.Range("L10").Value = 15000
$ws.Range("M10").Value = 14500
$ws.Range("P10").Value = 1115
$ws.Range("D11").Value = 44921
$ws.Range("J11").Value = 540
$ws.Range("K11").Value = 14000
$ws.Range("L11").Value = 15000
$ws.Range("M11").Value = 14500
$ws.Range("P11").Value = 1115
$ws.Range("D12").Value = 45077
$ws.Range("J12").Value = 360
$ws.Range("K12").Value = 17000
$ws.Range("L12").Value = 18000
$ws.Range("M12").Value = 17500
$ws.Range("P12").Value = 1346
$ws.Range("D13").Value = 44694
$ws.Range("J13").Value = 400
$ws.Range("K13").Value = 13000
$ws.Range("L13").Value = 14000
$ws.Range("M13").Value = 13500
$ws.Range("P13").Value = 1038
$ws.Range("D14").Value = 44756
$ws.Range("J14").Value = 400
$ws.Range("K14").Value = 14500
$ws.Range("L14").Value = 15000
$ws.Range("M14").Value = 14750
$ws.Range("P14").Value = 1135
$ws.Range("D15").Value = 44323
$ws.Range("J15").Value = 460
$ws.Range("K15").Value = 25000
$ws.Range("L15").Value = 26000
$ws.Range("M15").Value = 25500
$ws.Range("P15").Value = 1962
$ws.Range("D16").Value = 44757
$ws.Range("J16").Value = 400
$ws.Range("K16").Value = 15000
$ws.Range("L16").Value = 16000
$ws.Range("M16").Value = 15500
$ws.Range("P16").Value = 1192
$ws.Range("D17").Value = 44832
$ws.Range("J17").Value = 540
$ws.Range("K17").Value = 14000
$ws.Range("L17").Value = 15000
$ws.Range("M17").Value = 14500
$ws.Range("P17").Value = 1115
$ws.Range("D18").Value = 44418
$ws.Range("J18").Value = 500
$ws.Range("K18").Value = 14000
$ws.Range("L18").Value = 15000
$ws.Range("M18").Value = 14500
$ws.Range("P18").Value = 1115
$ws.Range("D19").Value = 44736
$ws.Range("J19").Value = 400
$ws.Range("K19").Value = 16000
$ws.Range("L19").Value = 17000
$ws.Range("M19").Value = 16500
$ws.Range("P19").Value = 1269
$ws.Range("D20").Value = 44998
$ws.Range("J20").Value = 440
$ws.Range("K20").Value = 17500
$ws.Range("L20").Value = 18000
$ws.Range("M20").Value = 17750
$ws.Range("P20").Value = 1365
$ws.Range("D21").Value = 45034
$ws.Range("J21").Value = 400
$ws.Range("K21").Value = 17000
$ws.Range("L21").Value = 18000
$ws.Range("M21").Value = 17500
$ws.Range("P21").Value = 1346
$ws.Range("D22").Value = 44904
$ws.Range("J22").Value = 400
$ws.Range("K22").Value = 14000
$ws.Range("L22").Value = 15000
$ws.Range("M22").Value = 14500
$ws.Range("P22").Value = 1115
$ws.Range("D23").Value = 44428
$ws.Range("J23").Value = 480
$ws.Range("K23").Value = 14000
$ws.Range("L23").Value = 15000
$ws.Range("M23").Value = 14500
$ws.Range("P23").Value = 1115
$ws.Range("D24").Value = 44670
$ws.Range("J24").Value = 480
$ws.Range("K24").Value = 14500
$ws.Range("L24").Value = 15000
$ws.Range("M24").Value = 14750
$ws.Range("P24").Value = 1135
$ws.Range("D25").Value = 44917
$ws.Range("J25").Value = 540
$ws.Range("K25").Value = 14000
$ws.Range("L25").Value = 15000
$ws.Range("M25").Value = 14500
$ws.Range("P25").Value = 1115
$ws.Range("D26").Value = 44873
$ws.Range("J26").Value = 480
$ws.Range("K26").Value = 14000
$ws.Range("L26").Value = 15000
$ws.Range("M26").Value = 14500
$ws.Range("P26").Value = 1115
$ws.Range("D27").Value = 44798
$ws.Range("J27").Value = 400
$ws.Range("K27").Value = 14000
$ws.Range("L27").Value = 15000
$ws.Range("M27").Value = 14500
$ws.Range("P27").Value = 1115
$ws.Range("D28").Value = 44335
$ws.Range("J28").Value = 480
$ws.Range("K28").Value = 24500
$ws.Range("L28").Value = 25000
$ws.Range("M28").Value = 24750
$ws.Range("P28").Value = 1904
$ws.Range("D29").Value = 44897
$ws.Range("J29").Value = 400
$ws.Range("K29").Value = 14000
$ws.Range("L29").Value = 15000
$ws.Range("M29").Value = 14500
$ws.Range("P29").Value = 1115
$ws.Range("D30").Value = 44943
$ws.Range("J30").Value = 400
$ws.Range("K30").Value = 14000
$ws.Range("L30").Value = 15000
$ws.Range("M30").Value = 14500
$ws.Range("P30").Value = 1115
$ws.Range("D31").Value = 44383
$ws.Range("J31").Value = 200
$ws.Range("K31").Value = 17000
$ws.Range("L31").Value = 18000
$ws.Range("M31").Value = 17500
$ws.Range("P31").Value = 1346
$ws.Range("D32").Value = 45068
$ws.Range("J32").Value = 460
$ws.Range("K32").Value = 17000
$ws.Range("L32").Value = 18000
$ws.Range("M32").Value = 17500
$ws.Range("P32").Value = 1346
$ws.Range("D33").Value = 45016
$ws.Range("J33").Value = 430
$ws.Range("K33").Value = 16000
$ws.Range("L33").Value = 17000
$ws.Range("M33").Value = 16500
$ws.Range("P33").Value = 1269
$ws.Range("D34").Value = 44412
$ws.Range("J34").Value = 600
$ws.Range("K34").Value = 14000
$ws.Range("L34").Value = 15000
$ws.Range("M34").Value = 14500
$ws.Range("P34").Value = 1115
$ws.Range("D35").Value = 45061
$ws.Range("J35").Value = 400
$ws.Range("K35").Value = 17000
$ws.Range("L35").Value = 18000
$ws.Range("M35").Value = 17500
$ws.Range("P35").Value = 1346
$ws.Range("D36").Value = 44761
$ws.Range("J36").Value = 480
$ws.Range("K36").Value = 14500
$ws.Range("L36").Value = 15000
$ws.Range("M36").Value = 14750
$ws.Range("P36").Value = 1135
$ws.Range("D37").Value = 44915
$ws.Range("J37").Value = 440
$ws.Range("K37").Value = 14000
$ws.Range("L37").Value = 15000
$ws.Range("M37").Value = 14500
$ws.Range("P37").Value = 1115
$ws.Range("D38").Value = 44442
$ws.Range("J38").Value = 460
$ws.Range("K38").Value = 14000
$ws.Range("L38").Value = 15000
$ws.Range("M38").Value = 14500
$ws.Range("P38").Value = 1115
$ws.Range("D39").Value = 45014
$ws.Range("J39").Value = 360
$ws.Range("K39").Value = 16000
$ws.Range("L39").Value = 17000
$ws.Range("M39").Value = 16500
$ws.Range("P39").Value = 1269
$ws.Range("D40").Value = 44879
$ws.Range("J40").Value = 400
$ws.Range("K40").Value = 13500
$ws.Range("L40").Value = 14000
$ws.Range("M40").Value = 13750
$ws.Range("P40").Value = 1058
$ws.Range("D41").Value = 44837
$ws.Range("J41").Value = 600
$ws.Range("K41").Value = 14000
$ws.Range("L41").Value = 15000
$ws.Range("M41").Value = 14500
$ws.Range("P41").Value = 1115
$ws.Range("D42").Value = 45082
$ws.Range("J42").Value = 440
$ws.Range("K42").Value = 17000
$ws.Range("L42").Value = 18000
$ws.Range("M42").Value = 17500
$ws.Range("P42").Value = 1346
$ws.Range("D43").Value = 44400
$ws.Range("J43").Value = 600
$ws.Range("K43").Value = 15000
$ws.Range("L43").Value = 16000
$ws.Range("M43").Value = 15500
$ws.Range("P43").Value = 1192
$ws.Range("D44").Value = 44426
$ws.Range("J44").Value = 460
$ws.Range("K44").Value = 14000
$ws.Range("L44").Value = 15000
$ws.Range("M44").Value = 14500
$ws.Range("P44").Value = 1115
$ws.Range("D45").Value = 44631
$ws.Range("J45").Value = 400
$ws.Range("K45").Value = 16000
$ws.Range("L45").Value = 17000
$ws.Range("M45").Value = 16500
$ws.Range("P45").Value = 1269
$ws.Range("D46").Value = 44309
$ws.Range("J46").Value = 400
$ws.Range("K46").Value = 26000
$ws.Range("L46").Value = 27000
$ws.Range("M46").Value = 26500
$ws.Range("P46").Value = 2038
$ws.Range("D47").Value = 45002
$ws.Range("J47").Value = 400
$ws.Range("K47").Value = 17000
$ws.Range("L47").Value = 18000
$ws.Range("M47").Value = 17500
$ws.Range("P47").Value = 1346
$ws.Range("D48").Value = 44445
$ws.Range("J48").Value = 600
$ws.Range("K48").Value = 13000
$ws.Range("L48").Value = 14000
$ws.Range("M48").Value = 13500
$ws.Range("P48").Value = 1038
$ws.Range("D49").Value = 45049
$ws.Range("J49").Value = 500
$ws.Range("K49").Value = 17000
$ws.Range("L49").Value = 18000
$ws.Range("M49").Value = 17500
$ws.Range("P49").Value = 1346
$ws.Range("D50").Value = 44771
$ws.Range("J50").Value = 480
$ws.Range("K50").Value = 14000
$ws.Range("L50").Value = 15000
$ws.Range("M50").Value = 14500
$ws.Range("P50").Value = 1115
$ws.Range("D51").Value = 44747
$ws.Range("J51").Value = 440
$ws.Range("K51").Value = 15000
$ws.Range("L51").Value = 16000
$ws.Range("M51").Value = 15500
$ws.Range("P51").Value = 1192
$ws.Range("D52").Value = 44855
$ws.Range("J52").Value = 500
$ws.Range("K52").Value = 13800
$ws.Range("L52").Value = 14000
$ws.Range("M52").Value = 13900
$ws.Range("P52").Value = 1069
$ws.Range("D53").Value = 44750
$ws.Range("J53").Value = 480
$ws.Range("K53").Value = 15000
$ws.Range("L53").Value = 16000
$ws.Range("M53").Value = 15500
$ws.Range("P53").Value = 1192
$ws.Range("D54").Value = 45079
$ws.Range("J54").Value = 400
$ws.Range("K54").Value = 17000
$ws.Range("L54").Value = 18000
$ws.Range("M54").Value = 17500
$ws.Range("P54").Value = 1346
$ws.Range("D55").Value = 44692
$ws.Range("J55").Value = 400
$ws.Range("K55").Value = 14000
$ws.Range("L55").Value = 15000
$ws.Range("M55").Value = 14500
$ws.Range("P55").Value = 1115
$ws.Range("D56").Value = 44767
$ws.Range("J56").Value = 600
$ws.Range("K56").Value = 15000
$ws.Range("L56").Value = 16000
$ws.Range("M56").Value = 15500
$ws.Range("P56").Value = 1192
$ws.Range("D57").Value = 44365
$ws.Range("J57").Value = 500
$ws.Range("K57").Value = 19500
$ws.Range("L57").Value = 20000
$ws.Range("M57").Value = 19750
$ws.Range("P57").Value = 1519
$ws.Range("D58").Value = 44988
$ws.Range("J58").Value = 400
$ws.Range("K58").Value = 19000
$ws.Range("L58").Value = 20000
$ws.Range("M58").Value = 19500
$ws.Range("P58").Value = 1500
$ws.Range("D59").Value = 44715
$ws.Range("J59").Value = 500
$ws.Range("K59").Value = 15000
$ws.Range("L59").Value = 16000
$ws.Range("M59").Value = 15500
$ws.Range("P59").Value = 1192
$ws.Range("D60").Value = 44925
$ws.Range("J60").Value = 400
$ws.Range("K60").Value = 14000
$ws.Range("L60").Value = 15000
$ws.Range("M60").Value = 14500
$ws.Range("P60").Value = 1115
$ws.Range("D61").Value = 44533
$ws.Range("J61").Value = 520
$ws.Range("K61").Value = 17000
$ws.Range("L61").Value = 18000
$ws.Range("M61").Value = 17500
$ws.Range("P61").Value = 1346
$ws.Range("D62").Value = 44599
$ws.Range("J62").Value = 400
$ws.Range("K62").Value = 15000
$ws.Range("L62").Value = 16000
$ws.Range("M62").Value = 15500
$ws.Range("P62").Value = 1192
$ws.Range("D63").Value = 44312
$ws.Range("J63").Value = 400
$ws.Range("K63").Value = 26000
$ws.Range("L63").Value = 27000
$ws.Range("M63").Value = 26500
$ws.Range("P63").Value = 2038
$ws.Range("D64").Value = 44924
$ws.Range("J64").Value = 480
$ws.Range("K64").Value = 14000
$ws.Range("L64").Value = 15000
$ws.Range("M64").Value = 14500
$ws.Range("P64").Value = 1115
$ws.Range("D65").Value = 44839
$ws.Range("J65").Value = 520
$ws.Range("K65").Value = 14000
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = 14500
$ws.Range("P65").Value = 1115
$ws.Range("D66").Value = 44984
$ws.Range("J66").Value = 400
$ws.Range("K66").Value = 18000
$ws.Range("L66").Value = 19000
$ws.Range("M66").Value = 18500
$ws.Range("P66").Value = 1423
$ws.Range("D67").Value = 44419
$ws.Range("J67").Value = 600
$ws.Range("K67").Value = 14000
$ws.Range("L67").Value = 15000
$ws.Range("M67").Value = 14500
$ws.Range("P67").Value = 1115
$ws.Range("D68").Value = 44847
$ws.Range("J68").Value = 400
$ws.Range("K68").Value = 13500
$ws.Range("L68").Value = 14000
$ws.Range("M68").Value = 13750
$ws.Range("P68").Value = 1058
$ws.Range("D69").Value = 44680
$ws.Range("J69").Value = 400
$ws.Range("K69").Value = 13500
$ws.Range("L69").Value = 14000
$ws.Range("M69").Value = 13750
$ws.Range("P69").Value = 1058
$ws.Range("D70").Value = 45007
$ws.Range("J70").Value = 200
$ws.Range("K70").Value = 16500
$ws.Range("L70").Value = 17000
$ws.Range("M70").Value = 16750
$ws.Range("P70").Value = 1288
$ws.Range("D71").Value = 44876
$ws.Range("J71").Value = 400
$ws.Range("K71").Value = 14000
$ws.Range("L71").Value = 15000
$ws.Range("M71").Value = 14500
$ws.Range("P71").Value = 1115
$ws.Range("D72").Value = 44704
$ws.Range("J72").Value = 400
$ws.Range("K72").Value = 13000
$ws.Range("L72").Value = 14000
$ws.Range("M72").Value = 13500
$ws.Range("P72").Value = 1038
$ws.Range("D73").Value = 44945
$ws.Range("J73").Value = 480
$ws.Range("K73").Value = 13500
$ws.Range("L73").Value = 14000
$ws.Range("M73").Value = 13750
$ws.Range("P73").Value = 1058
$ws.Range("D74").Value = 45040
$ws.Range("J74").Value = 200
$ws.Range("K74").Value = 17000
$ws.Range("L74").Value = 18000
$ws.Range("M74").Value = 17500
$ws.Range("P74").Value = 1346
$ws.Range("D75").Value = 44937
$ws.Range("J75").Value = 460
$ws.Range("K75").Value = 14000
$ws.Range("L75").Value = 15000
$ws.Range("M75").Value = 14500
$ws.Range("P75").Value = 1115
$ws.Range("D76").Value = 44918
$ws.Range("J76").Value = 450
$ws.Range("K76").Value = 14000
$ws.Range("L76").Value = 15000
$ws.Range("M76").Value = 14500
$ws.Range("P76").Value = 1115
$ws.Range("D77").Value = 44914
$ws.Range("J77").Value = 460
$ws.Range("K77").Value = 14000
$ws.Range("L77").Value = 15000
$ws.Range("M77").Value = 14500
$ws.Range("P77").Value = 1115
$ws.Range("D78").Value = 45035
$ws.Range("J78").Value = 400
$ws.Range("K78").Value = 17000
$ws.Range("L78").Value = 18000
$ws.Range("M78").Value = 17500
$ws.Range("P78").Value = 1346
$ws.Range("D79").Value = 44923
$ws.Range("J79").Value = 450
$ws.Range("K79").Value = 14000
$ws.Range("L79").Value = 15000
$ws.Range("M79").Value = 14500
$ws.Range("P79").Value = 1115
$ws.Range("D80").Value = 44923
$ws.Range("J80").Value = 450
$ws.Range("K80").Value = 14000
$ws.Range("L80").Value = 15000
$ws.Range("M80").Value = 14500
$ws.Range("P80").Value = 1115
$ws.Range("D81").Value = 44754
$ws.Range("J81").Value = 480
$ws.Range("K81").Value = 15000
$ws.Range("L81").Value = 16000
$ws.Range("M81").Value = 15500
$ws.Range("P81").Value = 1192
$ws.Range("D82").Value = 44963
$ws.Range("J82").Value = 400
$ws.Range("K82").Value = 23000
$ws.Range("L82").Value = 24000
$ws.Range("M82").Value = 23500
$ws.Range("P82").Value = 1808
$ws.Range("D83").Value = 44803
$ws.Range("J83").Value = 520
$ws.Range("K83").Value = 14000
$ws.Range("L83").Value = 15000
$ws.Range("M83").Value = 14500
$ws.Range("P83").Value = 1115
$ws.Range("D84").Value = 44763
$ws.Range("J84").Value = 500
$ws.Range("K84").Value = 15000
$ws.Range("L84").Value = 16000
$ws.Range("M84").Value = 15500
$ws.Range("P84").Value = 1192
$ws.Range("D85").Value = 44868
$ws.Range("J85").Value = 500
$ws.Range("K85").Value = 13500
$ws.Range("L85").Value = 14000
$ws.Range("M85").Value = 13750
$ws.Range("P85").Value = 1058
$ws.Range("D86").Value = 44826
$ws.Range("J86").Value = 520
$ws.Range("K86").Value = 14000
$ws.Range("L86").Value = 15000
$ws.Range("M86").Value = 14500
$ws.Range("P86").Value = 1115
$ws.Range("D87").Value = 44596
$ws.Range("J87").Value = 500
$ws.Range("K87").Value = 16000
$ws.Range("L87").Value = 17000
$ws.Range("M87").Value = 16500
$ws.Range("P87").Value = 1269
$ws.Range("D88").Value = 44946
$ws.Range("J88").Value = 400
$ws.Range("K88").Value = 13000
$ws.Range("L88").Value = 14000
$ws.Range("M88").Value = 13500
$ws.Range("P88").Value = 1038
$ws.Range("D89").Value = 44326
$ws.Range("J89").Value = 460
$ws.Range("K89").Value = 25000
$ws.Range("L89").Value = 26000
$ws.Range("M89").Value = 25500
$ws.Range("P89").Value = 1962
$ws.Range("D90").Value = 44740
$ws.Range("J90").Value = 500
$ws.Range("K90").Value = 16000
$ws.Range("L90").Value = 17000
$ws.Range("M90").Value = 16500
$ws.Range("P90").Value = 1269
$ws.Range("D91").Value = 44505
$ws.Range("J91").Value = 400
$ws.Range("K91").Value = 16000
$ws.Range("L91").Value = 17000
$ws.Range("M91").Value = 16500
$ws.Range("P91").Value = 1269
$ws.Range("D92").Value = 44810
$ws.Range("J92").Value = 540
$ws.Range("K92").Value = 14000
$ws.Range("L92").Value = 15000
$ws.Range("M92").Value = 14500
$ws.Range("P92").Value = 1115
$ws.Range("D93").Value = 45072
$ws.Range("J93").Value = 400
$ws.Range("K93").Value = 17000
$ws.Range("L93").Value = 18000
$ws.Range("M93").Value = 17500
$ws.Range("P93").Value = 1346
$ws.Range("D94").Value = 44860
$ws.Range("J94").Value = 400
$ws.Range("K94").Value = 14000
$ws.Range("L94").Value = 15000
$ws.Range("M94").Value = 14500
$ws.Range("P94").Value = 1115
$ws.Range("D95").Value = 44746
$ws.Range("J95").Value = 480
$ws.Range("K95").Value = 15000
$ws.Range("L95").Value = 16000
$ws.Range("M95").Value = 15500
$ws.Range("P95").Value = 1192
$ws.Range("D96").Value = 44708
$ws.Range("J96").Value = 440
$ws.Range("K96").Value = 13000
$ws.Range("L96").Value = 14000
$ws.Range("M96").Value = 13500
$ws.Range("P96").Value = 1038
$ws.Range("D97").Value = 44936
$ws.Range("J97").Value = 440
$ws.Range("K97").Value = 14000
$ws.Range("L97").Value = 15000
$ws.Range("M97").Value = 14500
$ws.Range("P97").Value = 1115
$ws.Range("D98").Value = 45051
$ws.Range("J98").Value = 340
$ws.Range("K98").Value = 17000
$ws.Range("L98").Value = 18000
$ws.Range("M98").Value = 17500
$ws.Range("P98").Value = 1346
$ws.Range("D99").Value = 44972
$ws.Range("J99").Value = 300
$ws.Range("K99").Value = 21000
$ws.Range("L99").Value = 22000
$ws.Range("M99").Value = 21500
$ws.Range("P99").Value = 1654
$ws.Range("D100").Value = 45070
$ws.Range("J100").Value = 320
$ws.Range("K100").Value = 17000
$ws.Range("L100").Value = 18000
$ws.Range("M100").Value = 17500
$ws.Range("P100").Value = 1346
$ws.Range("D101").Value = 45042
$ws.Range("J101").Value = 340
$ws.Range("K101").Value = 16000
$ws.Range("L101").Value = 17000
$ws.Range("M101").Value = 16500
$ws.Range("P101").Value = 1269
$ws.Range("D102").Value = 44582
$ws.Range("J102").Value = 520
$ws.Range("K102").Value = 15000
$ws.Range("L102").Value = 16000
$ws.Range("M102").Value = 15500
$ws.Range("P102").Value = 1192
$ws.Range("D103").Value = 45012
$ws.Range("J103").Value = 500
$ws.Range("K103").Value = 16000
$ws.Range("L103").Value = 17000
$ws.Range("M103").Value = 16500
$ws.Range("P103").Value = 1269
$ws.Range("D104").Value = 44657
$ws.Range("J104").Value = 460
$ws.Range("K104").Value = 15000
$ws.Range("L104").Value = 16000
$ws.Range("M104").Value = 15500
$ws.Range("P104").Value = 1192
$ws.Range("D105").Value = 44831
$ws.Range("J105").Value = 600
$ws.Range("K105").Value = 14000
$ws.Range("L105").Value = 15000
$ws.Range("M105").Value = 14500
$ws.Range("P105").Value = 1115
$ws.Range("D106").Value = 44687
$ws.Range("J106").Value = 440
$ws.Range("K106").Value = 14000
$ws.Range("L106").Value = 15000
$ws.Range("M106").Value = 14500
$ws.Range("P106").Value = 1115
$ws.Range("D107").Value = 44883
$ws.Range("J107").Value = 480
$ws.Range("K107").Value = 13500
$ws.Range("L107").Value = 14000
$ws.Range("M107").Value = 13750
$ws.Range("P107").Value = 1058
$ws.Range("D108").Value = 44939
$ws.Range("J108").Value = 460
$ws.Range("K108").Value = 14000
$ws.Range("L108").Value = 15000
$ws.Range("M108").Value = 14500
$ws.Range("P108").Value = 1115
$ws.Range("D109").Value = 45005
$ws.Range("J109").Value = 240
$ws.Range("K109").Value = 16000
$ws.Range("L109").Value = 17000
$ws.Range("M109").Value = 16500
$ws.Range("P109").Value = 1269
$ws.Range("D110").Value = 44957
$ws.Range("J110").Value = 400
$ws.Range("K110").Value = 23000
$ws.Range("L110").Value = 24000
$ws.Range("M110").Value = 23500
$ws.Range("P110").Value = 1808
$ws.Range("D111").Value = 44701
$ws.Range("J111").Value = 440
$ws.Range("K111").Value = 14000
$ws.Range("L111").Value = 15000
$ws.Range("M111").Value = 14500
$ws.Range("P111").Value = 1115
$ws.Range("D112").Value = 45033
$ws.Range("J112").Value = 400
$ws.Range("K112").Value = 17000
$ws.Range("L112").Value = 18000
$ws.Range("M112").Value = 17500
$ws.Range("P112").Value = 1346
$ws.Range("D113").Value = 44964
$ws.Range("J113").Value = 300
$ws.Range("K113").Value = 23000
$ws.Range("L113").Value = 24000
$ws.Range("M113").Value = 23500
$ws.Range("P113").Value = 1808
$ws.Range("D114").Value = 44910
$ws.Range("J114").Value = 460
$ws.Range("K114").Value = 14000
$ws.Range("L114").Value = 15000
$ws.Range("M114").Value = 14500
$ws.Range("P114").Value = 1115
$ws.Range("D115").Value = 44333
$ws.Range("J115").Value = 440
$ws.Range("K115").Value = 24000
$ws.Range("L115").Value = 25000
$ws.Range("M115").Value = 24500
$ws.Range("P115").Value = 1885
$ws.Range("D116").Value = 44875
$ws.Range("J116").Value = 480
$ws.Range("K116").Value = 14000
$ws.Range("L116").Value = 15000
$ws.Range("M116").Value = 14500
$ws.Range("P116").Value = 1115
$ws.Range("D117").Value = 44799
$ws.Range("J117").Value = 460
$ws.Range("K117").Value = 14000
$ws.Range("L117").Value = 15000
$ws.Range("M117").Value = 14500
$ws.Range("P117").Value = 1115
$ws.Range("D118").Value = 44344
$ws.Range("J118").Value = 400
$ws.Range("K118").Value = 18500
$ws.Range("L118").Value = 19000
$ws.Range("M118").Value = 18750
$ws.Range("P118").Value = 1442
$ws.Range("D119").Value = 44841
$ws.Range("J119").Value = 440
$ws.Range("K119").Value = 13500
$ws.Range("L119").Value = 14000
$ws.Range("M119").Value = 13750
$ws.Range("P119").Value = 1058
$ws.Range("D120").Value = 44484
$ws.Range("J120").Value = 360
$ws.Range("K120").Value = 14000
$ws.Range("L120").Value = 15000
$ws.Range("M120").Value = 14500
$ws.Range("P120").Value = 1115
$ws.Range("D121").Value = 45076
$ws.Range("J121").Value = 360
$ws.Range("K121").Value = 17000
$ws.Range("L121").Value = 18000
$ws.Range("M121").Value = 17500
$ws.Range("P121").Value = 1346
$ws.Range("D122").Value = 45021
$ws.Range("J122").Value = 400
$ws.Range("K122").Value = 17000
$ws.Range("L122").Value = 18000
$ws.Range("M122").Value = 17500
$ws.Range("P122").Value = 1346
$ws.Range("D123").Value = 44764
$ws.Range("J123").Value = 400
$ws.Range("K123").Value = 15000
$ws.Range("L123").Value = 16000
$ws.Range("M123").Value = 15500
$ws.Range("P123").Value = 1192
$ws.Range("D124").Value = 44414
$ws.Range("J124").Value = 500
$ws.Range("K124").Value = 14000
$ws.Range("L124").Value = 15000
$ws.Range("M124").Value = 14500
$ws.Range("P124").Value = 1115
$ws.Range("D125").Value = 44922
$ws.Range("J125").Value = 400
$ws.Range("K125").Value = 14000
$ws.Range("L125").Value = 15000
$ws.Range("M125").Value = 14500
$ws.Range("P125").Value = 1115
$ws.Range("D126").Value = 44846
$ws.Range("J126").Value = 400
$ws.Range("K126").Value = 13500
$ws.Range("L126").Value = 14000
$ws.Range("M126").Value = 13750
$ws.Range("P126").Value = 1058
$ws.Range("D127").Value = 44965
$ws.Range("J127").Value = 400
$ws.Range("K127").Value = 22500
$ws.Range("L127").Value = 23000
$ws.Range("M127").Value = 22750
$ws.Range("P127").Value = 1750
$ws.Range("D128").Value = 44379
$ws.Range("J128").Value = 600
$ws.Range("K128").Value = 17000
$ws.Range("L128").Value = 18000
$ws.Range("M128").Value = 17500
$ws.Range("P128").Value = 1346
$ws.Range("D129").Value = 44908
$ws.Range("J129").Value = 520
$ws.Range("K129").Value = 14000
$ws.Range("L129").Value = 15000
$ws.Range("M129").Value = 14500
$ws.Range("P129").Value = 1115
$ws.Range("D130").Value = 44792
$ws.Range("J130").Value = 400
$ws.Range("K130").Value = 14000
$ws.Range("L130").Value = 15000
$ws.Range("M130").Value = 14500
$ws.Range("P130").Value = 1115
$ws.Range("D131").Value = 44938
$ws.Range("J131").Value = 440
$ws.Range("K131").Value = 14000
$ws.Range("L131").Value = 15000
$ws.Range("M131").Value = 14500
$ws.Range("P131").Value = 1115
$ws.Range("D132").Value = 44644
$ws.Range("J132").Value = 400
$ws.Range("K132").Value = 15000
$ws.Range("L132").Value = 16000
$ws.Range("M132").Value = 15500
$ws.Range("P132").Value = 1192
$ws.Range("D133").Value = 44575
$ws.Range("J133").Value = 500
$ws.Range("K133").Value = 14000
$ws.Range("L133").Value = 15000
$ws.Range("M133").Value = 14500
$ws.Range("P133").Value = 1115
$ws.Range("D134").Value = 44981
$ws.Range("J134").Value = 400
$ws.Range("K134").Value = 19000
$ws.Range("L134").Value = 20000
$ws.Range("M134").Value = 19500
$ws.Range("P134").Value = 1500
$ws.Range("D135").Value = 45037
$ws.Range("J135").Value = 440
$ws.Range("K135").Value = 16000
$ws.Range("L135").Value = 17000
$ws.Range("M135").Value = 16500
$ws.Range("P135").Value = 1269
$ws.Range("D136").Value = 44942
$ws.Range("J136").Value = 440
$ws.Range("K136").Value = 14000
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = 14500
$ws.Range("P136").Value = 1115
$ws.Range("D137").Value = 44753
$ws.Range("J137").Value = 400
$ws.Range("K137").Value = 14500
$ws.Range("L137").Value = 15000
$ws.Range("M137").Value = 14750
$ws.Range("P137").Value = 1135
